# fix Num of Cluster(100 data)
#
# 1) "_1_3" sheet: cell BD36 gets a new descriptive marker string
#    (was the literal text "0").
# 2) "_1_4" sheet: a brand-new row 13 is appended, mirroring the
#    "marker" rows already present on the sheet (all "0" text markers
#    except for a single "1" in column BE, plus a descriptive note in
#    column A).

$wb   = $excel.ActiveWorkbook
$ws3  = $wb.Worksheets.Item("_1_3")
$ws4  = $wb.Worksheets.Item("_1_4")

# ---------------------------------------------------------------------
# 1) _1_4 - add row 13.
# ---------------------------------------------------------------------
# Seed the whole B13:CX13 block from an existing "all text zero" style
# range so every cell lands as a literal text "0" (matching how the
# rest of the sheet stores its 0/1 markers) instead of Excel silently
# re-typing the numeric-looking string as a real number.
# (Do this BEFORE touching _1_3!BD36 below, since this source range
# includes that very cell.)
$ws3.Range("B36:CX36").Copy($ws4.Range("B13:CX13"))

# That source row has two stray "1" markers of its own (CE36/CQ36) -
# stomp them back down to text "0" using another known text-"0" cell
# as the copy source (keeps the cell's type as text, not a number).
$zeroSrc = $ws4.Range("D12")
$zeroSrc.Copy($ws4.Range("CE13"))
$zeroSrc.Copy($ws4.Range("CQ13"))

# BE13 is the row's actual "1" marker - copy from an existing text-"1"
# cell so it is stored the same way (text, not a number).
$oneSrc = $ws4.Range("B12")
$oneSrc.Copy($ws4.Range("BE13"))

$triangle = [char]0x25B3

# Column A gets the descriptive note. Assigning text containing a
# newline makes Excel auto-size the row height (customHeight), so we
# AutoFit right after to drop that override and keep the row's height
# the same as every other (un-customized) row on the sheet.
$a13 = ", exhaust gas back pressure of 6.0 kPa at 100% load" + [char]10 + $triangle + "2"
$ws4.Range("A13").Value = $a13
$ws4.Rows.Item(13).AutoFit()

# ---------------------------------------------------------------------
# 2) _1_3!BD36 - replace the "0" marker with the real cluster text.
# ---------------------------------------------------------------------
$bd36 = "(+)25 & (+)% & (+),70 & (+)% & (+)" + $triangle + "1 & (-)1/2 & (-)3/4"
$ws3.Range("BD36").Value = $bd36
